# Feat: update timeline worksheet
#
# 1. Sort the existing Timeline data (rows 2:25) ascending by Start_Date
#    (column B), keeping the header row in place.
# 2. Insert two new rows right after the first sorted row for two new
#    timeline entries: a contract photography gig (PHOTO) and a freelance
#    web-dev/SEO gig (WEBDEV).
# 3. Grow Table1 so it covers the two new rows.
# 4. Restore the selection to A5 (matches the saved selection in the file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Sort existing data by Start_Date ascending (A1:G25 incl. header) ---
$sortRange = $ws.Range("A1:G25")
$sortKey = $ws.Range("B1:B25")
$sortRange.Sort($sortKey, 1, $null, $null, 1, $null, $null, 1)

# --- 2. Insert two blank rows after row 2, then populate them ---
$ws.Rows("3:4").Insert()

$ws.Range("A3").Value = "PHOTO"
$ws.Range("B3").Value = 44197
$ws.Range("C3").Value = 44834
$ws.Range("D3").Value = "Contract Photographer"
$ws.Range("E3").Value = "Photographer"
$ws.Range("F3").Value = "Fifty Mil Studios/KIP"
$ws.Range("G3").Value = "Present"

$ws.Range("A4").Value = "WEBDEV"
$ws.Range("B4").Value = 43831
$ws.Range("C4").Value = 44834
$ws.Range("D4").Value = "Web Development & SEO"
$ws.Range("E4").Value = "Web Dev."
$ws.Range("F4").Value = "Freelance"
$ws.Range("G4").Value = "Present"

# --- 3. Expand Table1 to include the new rows (A1:J27) ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J27"))

# --- 4. Restore selection ---
$ws.Range("A5").Select()
